$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 678.375
$ws.Range("I19").Value = 666
$ws.Range("J19").Value = 685.8
$ws.Range("K19").Value = 666
$ws.Range("L19").Value = 685.8
$ws.Range("M19").Value = -491

$ws.Range("H28").Value = 943.94116
$ws.Range("I28").Value = 943.94116
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 943.94116
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -458.94116

$ws.Range("H55").Value = 731.6667
$ws.Range("I55").Value = 731.6667
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 731.6667
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -517.6667

$ws.Range("H62").Value = 7995.2856
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7995.2856
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 7995.2856
$ws.Range("N62").Value = -9243.285599999999
$ws.Range("M62").Value = ""

$ws.Range("H65").Value = 7995.2856
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7995.2856
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 39976.428
$ws.Range("N65").Value = -46216.428
$ws.Range("M65").Value = ""

$ws.Range("H113").Value = 4266
$ws.Range("I113").Value = 4110
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 4110
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -856
$ws.Range("N113").Value = -11008

$ws.Range("H133").Value = 71499
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 71499
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 71499
$ws.Range("N133").Value = -81619

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1145.7273
$ws.Range("I2").Value = 567
$ws.Range("J2").Value = 3750
$ws.Range("K2").Value = 567
$ws.Range("L2").Value = 3750
$ws.Range("M2").Value = -454

$ws.Range("H45").Value = 1499
$ws.Range("I45").Value = 1499
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1499
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1122

$ws.Range("H74").Value = 2296.68
$ws.Range("I74").Value = 1735.1578
$ws.Range("J74").Value = 4074.8333
$ws.Range("K74").Value = 1735.1578
$ws.Range("L74").Value = 4074.8333
$ws.Range("M74").Value = -861.1578
$ws.Range("N74").Value = -5822.8333

$ws.Range("H76").Value = 42995
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 42995
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 42995
$ws.Range("N76").Value = -43671

$ws.Range("H77").Value = 2296.68
$ws.Range("I77").Value = 1735.1578
$ws.Range("J77").Value = 4074.8333
$ws.Range("K77").Value = 8675.789000000001
$ws.Range("L77").Value = 20374.1665
$ws.Range("M77").Value = -4307.789000000001
$ws.Range("N77").Value = -29110.1665

$ws.Range("H79").Value = 42995
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 42995
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 42995
$ws.Range("N79").Value = -45335

$ws.Range("H116").Value = 1145.7273
$ws.Range("I116").Value = 567
$ws.Range("J116").Value = 3750
$ws.Range("K116").Value = 567
$ws.Range("L116").Value = 3750
$ws.Range("M116").Value = 1727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1145.7273
$ws.Range("I3").Value = 567
$ws.Range("J3").Value = 3750
$ws.Range("K3").Value = 567
$ws.Range("L3").Value = 3750
$ws.Range("M3").Value = -453

$ws.Range("H22").Value = 605.35
$ws.Range("I22").Value = 497
$ws.Range("J22").Value = 806.5714
$ws.Range("K22").Value = 497
$ws.Range("L22").Value = 806.5714
$ws.Range("M22").Value = -324
$ws.Range("N22").Value = -1152.5714

$ws.Range("H64").Value = 279.85715
$ws.Range("I64").Value = 99
$ws.Range("J64").Value = 352.2
$ws.Range("K64").Value = 99
$ws.Range("L64").Value = 352.2
$ws.Range("M64").Value = 126
$ws.Range("N64").Value = -802.2

$ws.Range("H67").Value = 279.85715
$ws.Range("I67").Value = 99
$ws.Range("J67").Value = 352.2
$ws.Range("K67").Value = 99
$ws.Range("L67").Value = 352.2
$ws.Range("M67").Value = 681
$ws.Range("N67").Value = -1912.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 32112.75
$ws.Range("I22").Value = 1149.6666
$ws.Range("J22").Value = 125002
$ws.Range("K22").Value = 1149.6666
$ws.Range("L22").Value = 125002
$ws.Range("M22").Value = -799.6666
$ws.Range("N22").Value = -125702

$ws.Range("H99").Value = 2194.8572
$ws.Range("I99").Value = 1484.5
$ws.Range("J99").Value = 2479
$ws.Range("K99").Value = 1484.5
$ws.Range("L99").Value = 2479
$ws.Range("M99").Value = 13.5
$ws.Range("N99").Value = -5475

$ws.Range("H107").Value = 848.9286
$ws.Range("I107").Value = 684.2857
$ws.Range("J107").Value = 1013.5714
$ws.Range("K107").Value = 684.2857
$ws.Range("L107").Value = 1013.5714
$ws.Range("M107").Value = 1235.7143

$ws.Range("H126").Value = 2194.8572
$ws.Range("I126").Value = 1484.5
$ws.Range("J126").Value = 2479
$ws.Range("K126").Value = 4453.5
$ws.Range("L126").Value = 7437
$ws.Range("M126").Value = -1983.5
$ws.Range("N126").Value = -12377

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 247.625
$ws.Range("I92").Value = 172.125
$ws.Range("J92").Value = 474.125
$ws.Range("K92").Value = 516.375
$ws.Range("L92").Value = 1422.375
$ws.Range("M92").Value = 731.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7769
$ws.Range("I70").Value = 7203.5
$ws.Range("J70").Value = 8900
$ws.Range("K70").Value = 7203.5
$ws.Range("L70").Value = 8900
$ws.Range("M70").Value = -6933.5
$ws.Range("N70").Value = -9440

$ws.Range("H73").Value = 7769
$ws.Range("I73").Value = 7203.5
$ws.Range("J73").Value = 8900
$ws.Range("K73").Value = 7203.5
$ws.Range("L73").Value = 8900
$ws.Range("M73").Value = -6267.5
$ws.Range("N73").Value = -10772

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""

$ws.Range("H102").Value = 1432.1333
$ws.Range("I102").Value = 1031.9166
$ws.Range("J102").Value = 3033
$ws.Range("K102").Value = 1031.9166
$ws.Range("L102").Value = 3033
$ws.Range("M102").Value = 590.0834

$ws.Range("H113").Value = 1823.3334
$ws.Range("I113").Value = 1735
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1735
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 435

$ws.Range("H132").Value = 1251.5555
$ws.Range("I132").Value = 1233
$ws.Range("J132").Value = 1400
$ws.Range("K132").Value = 3699
$ws.Range("L132").Value = 4200
$ws.Range("M132").Value = -1169

$ws.Range("H135").Value = 526315
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 526315
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 526315
$ws.Range("N135").Value = -536455

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 950
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 950
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 950
$ws.Range("N12").Value = -1290
$ws.Range("M12").Value = ""

$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = ""

$ws.Range("H27").Value = 900
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -793
$ws.Range("N27").Value = ""

$ws.Range("H61").Value = 3823.625
$ws.Range("I61").Value = 3681.5
$ws.Range("J61").Value = 4250
$ws.Range("K61").Value = 3681.5
$ws.Range("L61").Value = 4250
$ws.Range("M61").Value = -3479.5
$ws.Range("N61").Value = -4654

$ws.Range("H113").Value = 3823.625
$ws.Range("I113").Value = 3681.5
$ws.Range("J113").Value = 4250
$ws.Range("K113").Value = 3681.5
$ws.Range("L113").Value = 4250
$ws.Range("M113").Value = -1511.5
$ws.Range("N113").Value = -8590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1403.68
$ws.Range("I136").Value = 1378.625
$ws.Range("J136").Value = 2005
$ws.Range("K136").Value = 4135.875
$ws.Range("L136").Value = 6015
$ws.Range("M136").Value = -1585.875
